# PowerShell Word COM-interop script implementing the Portuguese (pt)
# translation edits described by the diff. All target strings are unique
# within their containing paragraph, and each Find is scoped to that
# paragraph's Range so the identical English text that still appears
# (unchanged) in the other language sections of this multi-language
# template is left untouched.

$d = $word.ActiveDocument

function Replace-InParagraph {
    param(
        [int]$Index,
        [string]$OldText,
        [string]$NewText
    )
    $p = $d.Paragraphs.Item($Index)
    $r = $p.Range
    $found = $r.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, $NewText, 2)
    if (-not $found) {
        Write-Output "WARNING: not found in paragraph $Index -> $OldText"
    }
}

# 1. "English" heading -> "Inglês"
Replace-InParagraph 3 "English" "Inglês"

# 2. Subject line text
Replace-InParagraph 13 ": Meet our team in [CITY] | [DATE] " ": Conheça a nossa equipa em [CITY] | [DATE]"

# 3. Seminar invitation heading
Replace-InParagraph 16 "You’re invited to our Deriv Partner Seminar" "Está convidado(a) para o nosso Seminário de Parceiros da Deriv"

# 4. Greeting
Replace-InParagraph 18 "Dear [PARTNER NAME], " "Olá  [PARTNER NAME], "

# 5. Intro paragraph
Replace-InParagraph 20 "We’re excited to let you know that the Deriv Affiliate team will be in [CITY] in [MONTH] to meet with you, our valued partners!" "Estamos felizes por informar que a equipa de Afiliados da Deriv estará em [CIDADE] em [MÊS] para se reunir com os nossos parceiros!"

# 6. Location note (second <w:t> inside the run, after the <w:br/>)
Replace-InParagraph 24 "Your country manager will inform you about the exact location by [DATE]" "O gestor do seu país irá informá-lo(a) sobre o local exato até ao dia [DATE]"

# 7. Seminar description
Replace-InParagraph 30 "In this one-day seminar, we’ll be providing technical and marketing support, offering the opportunity to network with other partners over a delicious lunch as well as listening to your feedback about our partnership programmes. This is your chance to get your voice heard, which will help us plan future efforts to support you better. " "Neste seminário de 1 dia, terá acesso a suporte técnico e de Marketing, a oportunidade de interagir com outros parceiros durante o almoço e partilhar as suas sugestões com a nossa equipa. Esta é a sua grande oportunidade de fazer com que a sua voz seja ouvida, para assim planearmos ações futuras capazes de proporcionar um suporte melhor. "

# 8. RSVP sentence (before the bold [DATE] run)
Replace-InParagraph 32 "Please RSVP by submitting the registration form by " "Por favor, confirme a sua presença enviando o formulário de registo até ao dia "

# 9. Closing sentence of the RSVP paragraph (after the bold [DATE] run)
Replace-InParagraph 32 ". Please note that attendance is confirmed on a first come, first served basis. We look forward to seeing you there!" ". Tenha em atenção que a participação será confirmada por ordem de chegada. Esperamos vê-lo em breve!"

# 10. Button text ("Send my details") — first occurrence (English section)
Replace-InParagraph 33 "Send my details" "Enviar os meus dados"

# 11. Contact paragraph, first run
Replace-InParagraph 36 "If you have any questions, please contact us via " "Em caso de dúvida, entre em contato connosco através do "

# 12. Contact paragraph, " or " between the two hyperlinks
Replace-InParagraph 36 " or " " ou pelo "

# 13. Contact paragraph, final run
Replace-InParagraph 36 ". / If you have any questions, please contact your country manager, [NAME], at [EMAIL ADDRESS] or [WHATSAPP NO] (WhatsApp). " " no nosso website. / Caso tenha alguma dúvida, contacte o gestor do seu país [NAME] em [EMAIL ADDRESS] ou  [WHATSAPP NO] através do (Whatsapp). "

# 14. Button text ("Send my details") — second occurrence (inside the French
#     section of this template, which the original diff also overwrites).
Replace-InParagraph 91 "Send my details" "Enviar os meus dados"

Write-Output "Done."
